$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.273.51'
$ws.Cells.Item(2, 5).Value = '  +0.81%  '

$ws.Cells.Item(3, 4).Value = '2.444.41'
$ws.Cells.Item(3, 5).Value = '  +0.12%  '

$origStyle = $ws.Cells.Item(4, 4).Style
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 4).Style = $origStyle
$ws.Cells.Item(4, 5).Value = '  -0.14%  '

$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '571.38'
$ws.Cells.Item(5, 4).Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  +0.63%  '

$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '146.33'
$ws.Cells.Item(6, 4).Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  +0.38%  '

$ws.Cells.Item(7, 5).Value = '  +0.05%  '

$ws.Cells.Item(8, 5).Value = '  +1.11%  '

$ws.Cells.Item(9, 4).Value = '2.440.32'
$ws.Cells.Item(9, 5).Value = '  -0.19%  '

$ws.Cells.Item(10, 5).Value = '  -0.58%  '

$ws.Cells.Item(11, 5).Value = '  +1.21%  '

$ws.Cells.Item(12, 5).Value = '  -1.08%  '

$ws.Cells.Item(13, 5).Value = '  -0.13%  '

$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '27.00'
$ws.Cells.Item(14, 4).Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  +0.22%  '

$origStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.0000178'
$ws.Cells.Item(15, 4).Style = $origStyle
$ws.Cells.Item(15, 5).Value = '  -1.74%  '

$ws.Cells.Item(16, 4).Value = '2.876.28'
$ws.Cells.Item(16, 5).Value = '  +2.74%  '

$ws.Cells.Item(17, 4).Value = '62.912.32'
$ws.Cells.Item(17, 5).Value = '  +0.51%  '

$ws.Cells.Item(18, 4).Value = '2.470.11'
$ws.Cells.Item(18, 5).Value = '  +1.28%  '

$ws.Cells.Item(19, 5).Value = '  +0.81%  '

$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '7.32'
$ws.Cells.Item(20, 4).Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  +5.47%  '

$origStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '327.36'
$ws.Cells.Item(21, 4).Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  +0.98%  '

$ws.Cells.Item(22, 5).Value = '  +0.53%  '

$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '2.09'
$ws.Cells.Item(23, 4).Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  +13.27%  '

$origStyle = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '1.00'
$ws.Cells.Item(24, 4).Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  -0.05%  '

$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '65.34'
$ws.Cells.Item(25, 4).Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  -2.88%  '

$origStyle = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '618.55'
$ws.Cells.Item(26, 4).Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  +5.70%  '

$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '9.01'
$ws.Cells.Item(27, 4).Style = $origStyle
$ws.Cells.Item(27, 5).Value = '  +4.36%  '

$ws.Cells.Item(28, 5).Value = '  +1.58%  '

$ws.Cells.Item(29, 4).Value = '2.560.84'
$ws.Cells.Item(29, 5).Value = '  +0.02%  '

$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.50'
$ws.Cells.Item(30, 4).Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  +3.84%  '

$ws.Cells.Item(32, 5).Value = '  -1.90%  '

$ws.Cells.Item(33, 5).Value = '  -4.31%  '

$ws.Cells.Item(34, 5).Value = '  +0.16%  '

$ws.Cells.Item(35, 5).Value = '  +7.62%  '

$ws.Cells.Item(36, 5).Value = '  -0.05%  '

$ws.Cells.Item(37, 5).Value = '  +0.13%  '

$origStyle = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '5.41'
$ws.Cells.Item(39, 4).Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  +0.13%  '

$origStyle = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '18.74'
$ws.Cells.Item(40, 4).Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  -0.33%  '

$ws.Cells.Item(41, 2).Value = 'Monero'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '146.18'
$ws.Cells.Item(41, 4).Style = $origStyle
$ws.Cells.Item(41, 5).Value = '  -1.28%  '

$ws.Cells.Item(42, 2).Value = 'dogwifhat'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '2.70'
$ws.Cells.Item(42, 4).Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  +10.49%  '

$ws.Cells.Item(43, 5).Value = '  -1.25%  '

$ws.Cells.Item(44, 5).Value = '  -0.21%  '

$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '41.89'
$ws.Cells.Item(45, 4).Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  +0.64%  '

$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '148.69'
$ws.Cells.Item(46, 4).Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  -0.07%  '

$ws.Cells.Item(47, 5).Value = '  +2.00%  '

$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '21.18'
$ws.Cells.Item(48, 4).Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  +2.99%  '

$ws.Cells.Item(49, 5).Value = '  -0.26%  '

$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.600'
$ws.Cells.Item(50, 4).Style = $origStyle

$ws.Cells.Item(51, 5).Value = '  +0.61%  '
